# Refresh the lrc2p Inhba-Acvr2b LR-pair sheet with newly recomputed TPM-based
# NATMI metrics. Columns A-D (cluster/ligand/receptor labels) are untouched;
# only the numeric metric columns E:T for data rows 2-10 change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns E..T in order:
# E  Ligand-expressing cells
# F  Ligand detection rate
# G  Ligand average expression value
# H  Ligand total expression value
# I  Ligand derived specificity of average expression value
# J  Ligand derived specificity of total expression value
# K  Receptor-expressing cells
# L  Receptor detection rate
# M  Receptor average expression value
# N  Receptor total expression value
# O  Receptor derived specificity of average expression value
# P  Receptor derived specificity of total expression value
# Q  Edge average expression weight
# R  Edge total expression weight
# S  Edge average expression derived specificity
# T  Edge total expression derived specificity
$rows = @(
    @{ Row = 2;  Values = @(1, 0.3333333333333333, 0.05045533333333333, 0.151366, 0.004442474524580737, 0.004442474524580737, 3, 1, 1.315861666666667, 3.947585, 0.2754050739440597, 0.2754050739440597, 0.06639223901222222, 0.5975301511100001, 0.00122348002493676, 0.00122348002493676) }
    @{ Row = 3;  Values = @(1, 0.3333333333333333, 0.05045533333333333, 0.151366, 0.004442474524580737, 0.004442474524580737, 3, 1, 1.452872333333333, 4.358617, 0.3040809095127364, 0.3040809095127364, 0.07330515786911111, 0.659746420822, 0.001350871693921672, 0.001350871693921672) }
    @{ Row = 4;  Values = @(1, 0.3333333333333333, 0.05045533333333333, 0.151366, 0.004442474524580737, 0.004442474524580737, 3, 1, 2.009179666666667, 6.027539, 0.4205140165432039, 0.4205140165432039, 0.1013738298082222, 0.912364468274, 0.001868122805722306, 0.001868122805722306) }
    @{ Row = 5;  Values = @(3, 1, 8.433639666666666, 25.300919, 0.7425623198471305, 0.7425623198471305, 3, 1, 1.315861666666667, 3.947585, 0.2754050739440597, 0.2754050739440597, 11.09750314784611, 99.87752833061501, 0.2045054306055715, 0.2045054306055715) }
    @{ Row = 6;  Values = @(3, 1, 8.433639666666666, 25.300919, 0.7425623198471305, 0.7425623198471305, 3, 1, 1.452872333333333, 4.358617, 0.3040809095127364, 0.3040809095127364, 12.25300174100255, 110.277015669023, 0.2257990255890029, 0.2257990255890029) }
    @{ Row = 7;  Values = @(3, 1, 8.433639666666666, 25.300919, 0.7425623198471305, 0.7425623198471305, 3, 1, 2.009179666666667, 6.027539, 0.4205140165432039, 0.4205140165432039, 16.94469733426011, 152.502276008341, 0.3122578636525561, 0.3122578636525561) }
    @{ Row = 8;  Values = @(3, 1, 2.873389, 8.620167, 0.2529952056282888, 0.2529952056282888, 3, 1, 1.315861666666667, 3.947585, 0.2754050739440597, 0.2754050739440597, 3.780982438521667, 34.028841946695, 0.06967616331355148, 0.06967616331355148) }
    @{ Row = 9;  Values = @(3, 1, 2.873389, 8.620167, 0.2529952056282888, 0.2529952056282888, 3, 1, 1.452872333333333, 4.358617, 0.3040809095127364, 0.3040809095127364, 4.174667381004333, 37.572006429039, 0.07693101222981183, 0.07693101222981183) }
    @{ Row = 10; Values = @(3, 1, 2.873389, 8.620167, 0.2529952056282888, 0.2529952056282888, 3, 1, 2.009179666666667, 6.027539, 0.4205140165432039, 0.4205140165432039, 5.773154753223667, 51.958392779013, 0.1063880300849255, 0.1063880300849255) }
)

# Columns E(5) through T(20), written one cell at a time for maximum
# COM-interop compatibility (avoids relying on 2-D SAFEARRAY marshalling).
foreach ($entry in $rows) {
    $r = $entry.Row
    $colIndex = 5
    foreach ($val in $entry.Values) {
        $ws.Cells.Item($r, $colIndex).Value2 = $val
        $colIndex++
    }
}
